# Rename the "Roster" worksheet to "01"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roster")
$ws.Name = "01"
